$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Delete the "Code" column (column G), which holds the 0x60 "Code" values.
# This shifts the "Operandos" column (and everything to its right) one
# column to the left, matching the edit described in the commit (removal
# of the "Code"/"Funct value" column from the instruction tables).
$ws.Columns("G").Delete() | Out-Null

# Update the selection to match the target state.
$ws.Range("E7").Select() | Out-Null
